# Starting to check-in changes for preprocessing of data
#
# - D10 no longer needs to be flagged red: reset its font/format to the
#   same "normal" black + wrapped + bordered look already used by D9.
# - Column F gets a "Status" value for the two "merging documents" rows
#   (11 and 12): "Not Needed" / "Noit Needed".
# - Column F is widened slightly so the new values are readable.
# - Selection left on D10 (last cell the author touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-flag D10 (was red/bold-ish warning text) by copying the plain
# black/bordered/wrapped formatting already used elsewhere in column D.
$ws.Range("D9").Copy()
$ws.Range("D10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New "Status" values for rows 11 and 12 (column F).
$ws.Range("F11").Value = "Not Needed"
$ws.Range("F12").Value = "Noit Needed"

# Give column F an explicit width now that it holds real content.
$ws.Columns("F").ColumnWidth = 12.1667

# Leave the selection where the author left it.
$ws.Range("D10").Select()
